$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-23 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-24 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("5+71=", $true, $false, $false, $false, $false, $true, 1, $false, "83-24=", 2) | Out-Null
$d.Content.Find.Execute("59-30=", $true, $false, $false, $false, $false, $true, 1, $false, "62-51=", 2) | Out-Null
$d.Content.Find.Execute("77-28=", $true, $false, $false, $false, $false, $true, 1, $false, "27+40=", 2) | Out-Null
$d.Content.Find.Execute("28+14=", $true, $false, $false, $false, $false, $true, 1, $false, "28+56=", 2) | Out-Null
$d.Content.Find.Execute("78-50=", $true, $false, $false, $false, $false, $true, 1, $false, "58-36=", 2) | Out-Null
$d.Content.Find.Execute("22-0=", $true, $false, $false, $false, $false, $true, 1, $false, "38+27=", 2) | Out-Null
$d.Content.Find.Execute("53+24=", $true, $false, $false, $false, $false, $true, 1, $false, "11+37=", 2) | Out-Null
$d.Content.Find.Execute("7+61=", $true, $false, $false, $false, $false, $true, 1, $false, "88-45=", 2) | Out-Null
$d.Content.Find.Execute("23+76=", $true, $false, $false, $false, $false, $true, 1, $false, "17-0=", 2) | Out-Null
$d.Content.Find.Execute("99-98=", $true, $false, $false, $false, $false, $true, 1, $false, "81-18=", 2) | Out-Null
$d.Content.Find.Execute("38+26=", $true, $false, $false, $false, $false, $true, 1, $false, "82-63=", 2) | Out-Null
$d.Content.Find.Execute("96-60=", $true, $false, $false, $false, $false, $true, 1, $false, "65-52=", 2) | Out-Null
$d.Content.Find.Execute("35-33=", $true, $false, $false, $false, $false, $true, 1, $false, "69+7=", 2) | Out-Null
$d.Content.Find.Execute("43+53=", $true, $false, $false, $false, $false, $true, 1, $false, "67-14=", 2) | Out-Null
$d.Content.Find.Execute("28+69=", $true, $false, $false, $false, $false, $true, 1, $false, "5+70=", 2) | Out-Null
$d.Content.Find.Execute("39+44=", $true, $false, $false, $false, $false, $true, 1, $false, "10+51=", 2) | Out-Null
$d.Content.Find.Execute("92-48=", $true, $false, $false, $false, $false, $true, 1, $false, "16-15=", 2) | Out-Null
$d.Content.Find.Execute("4+76=", $true, $false, $false, $false, $false, $true, 1, $false, "19+26=", 2) | Out-Null
$d.Content.Find.Execute("81-8=", $true, $false, $false, $false, $false, $true, 1, $false, "7+56=", 2) | Out-Null
$d.Content.Find.Execute("18+60=", $true, $false, $false, $false, $false, $true, 1, $false, "0+92=", 2) | Out-Null
$d.Content.Find.Execute("68-42=", $true, $false, $false, $false, $false, $true, 1, $false, "55+32=", 2) | Out-Null
$d.Content.Find.Execute("73-61=", $true, $false, $false, $false, $false, $true, 1, $false, "96+3=", 2) | Out-Null
$d.Content.Find.Execute("74+25=", $true, $false, $false, $false, $false, $true, 1, $false, "4+50=", 2) | Out-Null
$d.Content.Find.Execute("91-83=", $true, $false, $false, $false, $false, $true, 1, $false, "35-9=", 2) | Out-Null
$d.Content.Find.Execute("17+30=", $true, $false, $false, $false, $false, $true, 1, $false, "55-10=", 2) | Out-Null
$d.Content.Find.Execute("18+62=", $true, $false, $false, $false, $false, $true, 1, $false, "0+27=", 2) | Out-Null
$d.Content.Find.Execute("47-0=", $true, $false, $false, $false, $false, $true, 1, $false, "25+68=", 2) | Out-Null
$d.Content.Find.Execute("89-75=", $true, $false, $false, $false, $false, $true, 1, $false, "42+10=", 2) | Out-Null
$d.Content.Find.Execute("96-53=", $true, $false, $false, $false, $false, $true, 1, $false, "65+17=", 2) | Out-Null
$d.Content.Find.Execute("26+23=", $true, $false, $false, $false, $false, $true, 1, $false, "55+28=", 2) | Out-Null
$d.Content.Find.Execute("49+15=", $true, $false, $false, $false, $false, $true, 1, $false, "96-44=", 2) | Out-Null
$d.Content.Find.Execute("67-24=", $true, $false, $false, $false, $false, $true, 1, $false, "36-2=", 2) | Out-Null
$d.Content.Find.Execute("29-7=", $true, $false, $false, $false, $false, $true, 1, $false, "7+84=", 2) | Out-Null
$d.Content.Find.Execute("21-6=", $true, $false, $false, $false, $false, $true, 1, $false, "57-53=", 2) | Out-Null
$d.Content.Find.Execute("88-3=", $true, $false, $false, $false, $false, $true, 1, $false, "94-64=", 2) | Out-Null
$d.Content.Find.Execute("39+2=", $true, $false, $false, $false, $false, $true, 1, $false, "97-79=", 2) | Out-Null
$d.Content.Find.Execute("57+28=", $true, $false, $false, $false, $false, $true, 1, $false, "56-3=", 2) | Out-Null
$d.Content.Find.Execute("82+16=", $true, $false, $false, $false, $false, $true, 1, $false, "5-4=", 2) | Out-Null
$d.Content.Find.Execute("42+26=", $true, $false, $false, $false, $false, $true, 1, $false, "55+8=", 2) | Out-Null
$d.Content.Find.Execute("65-13=", $true, $false, $false, $false, $false, $true, 1, $false, "17+25=", 2) | Out-Null
$d.Content.Find.Execute("30+16=", $true, $false, $false, $false, $false, $true, 1, $false, "77-24=", 2) | Out-Null
$d.Content.Find.Execute("11+41=", $true, $false, $false, $false, $false, $true, 1, $false, "76-32=", 2) | Out-Null
$d.Content.Find.Execute("27+62=", $true, $false, $false, $false, $false, $true, 1, $false, "27+44=", 2) | Out-Null
$d.Content.Find.Execute("29+13=", $true, $false, $false, $false, $false, $true, 1, $false, "56-17=", 2) | Out-Null
$d.Content.Find.Execute("75-31=", $true, $false, $false, $false, $false, $true, 1, $false, "22+29=", 2) | Out-Null
$d.Content.Find.Execute("7+19=", $true, $false, $false, $false, $false, $true, 1, $false, "88-53=", 2) | Out-Null
$d.Content.Find.Execute("32-28=", $true, $false, $false, $false, $false, $true, 1, $false, "3+58=", 2) | Out-Null
$d.Content.Find.Execute("25-8=", $true, $false, $false, $false, $false, $true, 1, $false, "68+24=", 2) | Out-Null
$d.Content.Find.Execute("11+31=", $true, $false, $false, $false, $false, $true, 1, $false, "61-17=", 2) | Out-Null
$d.Content.Find.Execute("34+7=", $true, $false, $false, $false, $false, $true, 1, $false, "4+19=", 2) | Out-Null
$d.Content.Find.Execute("78-38=", $true, $false, $false, $false, $false, $true, 1, $false, "32+16=", 2) | Out-Null
$d.Content.Find.Execute("15+67=", $true, $false, $false, $false, $false, $true, 1, $false, "71-18=", 2) | Out-Null
$d.Content.Find.Execute("85-9=", $true, $false, $false, $false, $false, $true, 1, $false, "22-15=", 2) | Out-Null
$d.Content.Find.Execute("65+3=", $true, $false, $false, $false, $false, $true, 1, $false, "59-18=", 2) | Out-Null
$d.Content.Find.Execute("26+69=", $true, $false, $false, $false, $false, $true, 1, $false, "5+76=", 2) | Out-Null
$d.Content.Find.Execute("59-12=", $true, $false, $false, $false, $false, $true, 1, $false, "85-61=", 2) | Out-Null
$d.Content.Find.Execute("7+33=", $true, $false, $false, $false, $false, $true, 1, $false, "6+39=", 2) | Out-Null
$d.Content.Find.Execute("88-76=", $true, $false, $false, $false, $false, $true, 1, $false, "49+1=", 2) | Out-Null
$d.Content.Find.Execute("1+75=", $true, $false, $false, $false, $false, $true, 1, $false, "35-20=", 2) | Out-Null
$d.Content.Find.Execute("90-56=", $true, $false, $false, $false, $false, $true, 1, $false, "94-64=", 2) | Out-Null
$d.Content.Find.Execute("26+33=", $true, $false, $false, $false, $false, $true, 1, $false, "99-6=", 2) | Out-Null
$d.Content.Find.Execute("49-17=", $true, $false, $false, $false, $false, $true, 1, $false, "97-94=", 2) | Out-Null
$d.Content.Find.Execute("47-3=", $true, $false, $false, $false, $false, $true, 1, $false, "68-49=", 2) | Out-Null
$d.Content.Find.Execute("33+4=", $true, $false, $false, $false, $false, $true, 1, $false, "73+21=", 2) | Out-Null
$d.Content.Find.Execute("1+94=", $true, $false, $false, $false, $false, $true, 1, $false, "23-22=", 2) | Out-Null
$d.Content.Find.Execute("20+71=", $true, $false, $false, $false, $false, $true, 1, $false, "62-23=", 2) | Out-Null
$d.Content.Find.Execute("85-50=", $true, $false, $false, $false, $false, $true, 1, $false, "8+53=", 2) | Out-Null
$d.Content.Find.Execute("42-32=", $true, $false, $false, $false, $false, $true, 1, $false, "34+6=", 2) | Out-Null
$d.Content.Find.Execute("31-6=", $true, $false, $false, $false, $false, $true, 1, $false, "68+11=", 2) | Out-Null
$d.Content.Find.Execute("50-46=", $true, $false, $false, $false, $false, $true, 1, $false, "83-78=", 2) | Out-Null
$d.Content.Find.Execute("6-5=", $true, $false, $false, $false, $false, $true, 1, $false, "13+27=", 2) | Out-Null
$d.Content.Find.Execute("29+19=", $true, $false, $false, $false, $false, $true, 1, $false, "41+33=", 2) | Out-Null
$d.Content.Find.Execute("54+3=", $true, $false, $false, $false, $false, $true, 1, $false, "3+86=", 2) | Out-Null
$d.Content.Find.Execute("68-60=", $true, $false, $false, $false, $false, $true, 1, $false, "12+70=", 2) | Out-Null
$d.Content.Find.Execute("27+3=", $true, $false, $false, $false, $false, $true, 1, $false, "84-29=", 2) | Out-Null
$d.Content.Find.Execute("21-18=", $true, $false, $false, $false, $false, $true, 1, $false, "57+35=", 2) | Out-Null
$d.Content.Find.Execute("44+10=", $true, $false, $false, $false, $false, $true, 1, $false, "72-12=", 2) | Out-Null
$d.Content.Find.Execute("41+12=", $true, $false, $false, $false, $false, $true, 1, $false, "34+17=", 2) | Out-Null
$d.Content.Find.Execute("69-59=", $true, $false, $false, $false, $false, $true, 1, $false, "14-12=", 2) | Out-Null
$d.Content.Find.Execute("51-9=", $true, $false, $false, $false, $false, $true, 1, $false, "73-9=", 2) | Out-Null
$d.Content.Find.Execute("30+37=", $true, $false, $false, $false, $false, $true, 1, $false, "90-32=", 2) | Out-Null
$d.Content.Find.Execute("7+34=", $true, $false, $false, $false, $false, $true, 1, $false, "49-5=", 2) | Out-Null
$d.Content.Find.Execute("24+40=", $true, $false, $false, $false, $false, $true, 1, $false, "62+18=", 2) | Out-Null
$d.Content.Find.Execute("25+39=", $true, $false, $false, $false, $false, $true, 1, $false, "77-74=", 2) | Out-Null
$d.Content.Find.Execute("18+67=", $true, $false, $false, $false, $false, $true, 1, $false, "33-5=", 2) | Out-Null
$d.Content.Find.Execute("12+44=", $true, $false, $false, $false, $false, $true, 1, $false, "97-84=", 2) | Out-Null
$d.Content.Find.Execute("41+2=", $true, $false, $false, $false, $false, $true, 1, $false, "93-13=", 2) | Out-Null
$d.Content.Find.Execute("54-33=", $true, $false, $false, $false, $false, $true, 1, $false, "95-74=", 2) | Out-Null
$d.Content.Find.Execute("54-8=", $true, $false, $false, $false, $false, $true, 1, $false, "33+38=", 2) | Out-Null
$d.Content.Find.Execute("36+23=", $true, $false, $false, $false, $false, $true, 1, $false, "63+23=", 2) | Out-Null
$d.Content.Find.Execute("42-23=", $true, $false, $false, $false, $false, $true, 1, $false, "7+14=", 2) | Out-Null
$d.Content.Find.Execute("16+45=", $true, $false, $false, $false, $false, $true, 1, $false, "70-59=", 2) | Out-Null
$d.Content.Find.Execute("99-89=", $true, $false, $false, $false, $false, $true, 1, $false, "86-35=", 2) | Out-Null
$d.Content.Find.Execute("41+13=", $true, $false, $false, $false, $false, $true, 1, $false, "58-50=", 2) | Out-Null
$d.Content.Find.Execute("6-1=", $true, $false, $false, $false, $false, $true, 1, $false, "29+43=", 2) | Out-Null
$d.Content.Find.Execute("74+6=", $true, $false, $false, $false, $false, $true, 1, $false, "10+77=", 2) | Out-Null
$d.Content.Find.Execute("28+10=", $true, $false, $false, $false, $false, $true, 1, $false, "45-1=", 2) | Out-Null
$d.Content.Find.Execute("90-75=", $true, $false, $false, $false, $false, $true, 1, $false, "49+32=", 2) | Out-Null
$d.Content.Find.Execute("77-69=", $true, $false, $false, $false, $false, $true, 1, $false, "84-29=", 2) | Out-Null
$d.Content.Find.Execute("41-35=", $true, $false, $false, $false, $false, $true, 1, $false, "54-6=", 2) | Out-Null
